$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Stash the two existing header formats (Bold 20 Calibri on B2, Bold 15
# Arial centered on C2) onto holding cells before we overwrite/clear the
# original range, so PasteSpecial(Formats) can reuse the existing style
# records instead of minting duplicate fonts/xfs. ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("Z2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Wipe the old B2:C6 block entirely (values + formats) ---
$ws.Range("B2:D15").Clear()

# --- Column widths (A is new). The host's ColumnWidth setter quantizes to
# 1/6-character steps, so these inputs are chosen to land on the closest
# achievable stored width to the target (6.5703125 / 25.140625 / 51.5703125 / 36).
$ws.Columns.Item(1).ColumnWidth = 5.65
$ws.Columns.Item(2).ColumnWidth = 24.325
$ws.Columns.Item(3).ColumnWidth = 50.65
$ws.Columns.Item(4).ColumnWidth = 35.15

# --- Row 2 header values ---
$ws.Range("B2").Value = "JumpPark:"
$ws.Range("C2").Value = "Aplicativo Operacional"
$ws.Range("D2").Value = "Clientes"
$ws.Rows.Item(2).RowHeight = 23.25

# --- Column C body (Aplicativo Operacional) ---
$ws.Range("C3").Value = ">Controla a entrada e saida dos carros"
$ws.Range("C4").Value = ">Auto peenche dados de clientes frequentes"
$ws.Range("C5").Value = ">App para smartphone ou tablet"
$ws.Range("C6").Value = ">Realiza todas as transações nescessárias"
$ws.Range("C7").Value = ">todas são transacionadas com saldo positivo ou negativo"
$ws.Range("C8").Value = ">Comunicação por email"
$ws.Range("C9").Value = ">Aviso de entrada e saida"
$ws.Range("C10").Value = ">Recibo de pagamento"
$ws.Range("C11").Value = ">Inicio e fim dos serviços"
$ws.Range("C12").Value = ">Checklist de avaria dos veículos"
$ws.Range("C13").Value = ">Mapeamento do veículo na entrada utilizando fotos"
$ws.Range("C14").Value = ">Inclusão de serviços automotivos"
$ws.Range("C15").Value = ">Idependencia de internet durante o uso do programa"

# --- Column D body (Clientes) ---
$ws.Range("D3").Value = ">Controle de convênios ou pós-pago"
$ws.Range("D4").Value = ">Geração de relatórios detalhados"
$ws.Range("D5").Value = ">Envio de recibos das faturas por e-mail "
$ws.Range("D6").Value = ">"

# --- Re-apply the stashed header formats onto the (re-created) header cells ---
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("Z2").Copy() | Out-Null
$ws.Range("C2:D2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1:Z2").Clear()

# --- Hyperlink cell B3 (value first, then Add so no stray display text is written) ---
$ws.Range("B3").Value = "https://jumppark.com.br/"
$ws.Hyperlinks.Add($ws.Range("B3"), "https://jumppark.com.br/") | Out-Null

$ws.Range("B10").Select()
